# remove weird kid ideas (pear123) and replace with numeric (idXX)
#
# The transcript used placeholder "kid" names (apple333, pear444) both as a
# standalone value in the "speaker" column and embedded inside narration
# text such as " Hi, [apple333]." Replace every occurrence workbook-wide
# with the de-identified numeric ids (id90 / id91), matching substrings
# wherever they occur (not just whole-cell matches).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Cells.Replace("apple333", "id90")
[void]$ws.Cells.Replace("pear444", "id91")
